$d = $word.ActiveDocument

# Change 1: "(VSP ViSion)" -> "(VSP Vision)" (capital S to lowercase s)
$d.Content.Find.Execute("VSP ViSion)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "VSP Vision)", 2)

# Change 2: "the purpose \u201cfuture engagements" -> "the purpose of \u201cfuture engagements"
$d.Content.Find.Execute("the purpose " + [char]8220 + "future engagements", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the purpose of " + [char]8220 + "future engagements", 2)

# Change 3: "results \u201c" -> "results of \u201c"
$d.Content.Find.Execute("desired results " + [char]8220, $true, $false, $false, $false, $false,
                         $true, 1, $false, "desired results of " + [char]8220, 2)
